$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09024753013142529
$ws.Range("C2").Value = 0.9987781041500535
$ws.Range("D2").Value = 0.2392058410935043
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
$ws.Range("G2").Value = 0.124317388383497
$ws.Range("H2").Value = 0.991
